$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 3624
$ws.Range("E2").Value = 888
$ws.Range("F2").Value = 888
$ws.Range("G2").Value = 601
$ws.Range("H2").Value = 422
$ws.Range("I2").Value = 245
$ws.Range("J2").Value = 176
$ws.Range("K2").Value = 7664
$ws.Range("L2").Value = 1534
$ws.Range("M2").Value = 6130
$ws.Range("N2").Value = 5660
$ws.Range("O2").Value = 470
$ws.Range("P2").Value = 13
$ws.Range("Q2").Value = 693
$ws.Range("R2").Value = -4626
$ws.Range("S2").Value = 4534
$ws.Range("T2").Value = 74
$ws.Range("U2").Value = 619
$ws.Range("V2").Value = 136
$ws.Range("W2").Value = 24.5
$ws.Range("X2").Value = 11.64
$ws.Range("Y2").Value = 7.82
$ws.Range("Z2").Value = 8.82
$ws.Range("AA2").Value = 25.03
$ws.Range("AB2").Value = 40939.77
$ws.Range("AC2").Value = 536
$ws.Range("AE2").Value = 10536
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 53716456
$ws.Range("AD2").ClearContents()
$ws.Range("AH2").ClearContents()

# Row 3
$ws.Range("D3").Value = 10729
$ws.Range("E3").Value = 2254
$ws.Range("F3").Value = 2254
$ws.Range("G3").Value = 2289
$ws.Range("H3").Value = 1686
$ws.Range("I3").Value = 1205
$ws.Range("J3").Value = 481
$ws.Range("K3").Value = 14584
$ws.Range("L3").Value = 2777
$ws.Range("M3").Value = 11807
$ws.Range("N3").Value = 10592
$ws.Range("O3").Value = 1215
$ws.Range("P3").Value = 15
$ws.Range("Q3").Value = 2243
$ws.Range("R3").Value = -4784
$ws.Range("S3").Value = 3695
$ws.Range("T3").Value = 187
$ws.Range("U3").Value = 2055
$ws.Range("V3").Value = 312
$ws.Range("W3").Value = 21
$ws.Range("X3").Value = 15.72
$ws.Range("Y3").Value = 14.83
$ws.Range("Z3").Value = 15.16
$ws.Range("AA3").Value = 23.52
$ws.Range("AB3").Value = 70483.45
$ws.Range("AC3").Value = 2049
$ws.Range("AE3").Value = 17786
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 59552592
$ws.Range("AD3").ClearContents()
$ws.Range("AH3").ClearContents()

# Row 4
$ws.Range("D4").Value = 15000
$ws.Range("E4").Value = 2947
$ws.Range("F4").Value = 2947
$ws.Range("G4").Value = 2756
$ws.Range("H4").Value = 2092
$ws.Range("I4").Value = 1740
$ws.Range("J4").Value = 353
$ws.Range("K4").Value = 19574
$ws.Range("L4").Value = 6471
$ws.Range("M4").Value = 13104
$ws.Range("N4").Value = 12163
$ws.Range("O4").Value = 940
$ws.Range("P4").Value = 68
$ws.Range("Q4").Value = 1971
$ws.Range("R4").Value = -2073
$ws.Range("S4").Value = 630
$ws.Range("T4").Value = 898
$ws.Range("U4").Value = 1072
$ws.Range("V4").Value = 1053
$ws.Range("W4").Value = 19.64
$ws.Range("X4").Value = 13.95
$ws.Range("Y4").Value = 15.29
$ws.Range("Z4").Value = 12.25
$ws.Range("AA4").Value = 49.38
$ws.Range("AB4").Value = 22045.83
$ws.Range("AC4").Value = 2719
$ws.Range("AE4").Value = 17946
$ws.Range("AF4").Value = 0
$ws.Range("AG4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 67776388
$ws.Range("AD4").ClearContents()
$ws.Range("AH4").ClearContents()

# Row 5
$ws.Range("D5").Value = 24248
$ws.Range("E5").Value = 5098
$ws.Range("F5").Value = 5098
$ws.Range("G5").Value = 4715
$ws.Range("H5").Value = 3609
$ws.Range("I5").Value = 3098
$ws.Range("J5").Value = 511
$ws.Range("K5").Value = 53477
$ws.Range("L5").Value = 8824
$ws.Range("M5").Value = 44653
$ws.Range("N5").Value = 43295
$ws.Range("O5").Value = 1358
$ws.Range("P5").Value = 85
$ws.Range("Q5").Value = 5034
$ws.Range("R5").Value = -14327
$ws.Range("S5").Value = 25273
$ws.Range("T5").Value = 300
$ws.Range("U5").Value = 4734
$ws.Range("V5").Value = 39
$ws.Range("W5").Value = 21.02
$ws.Range("X5").Value = 14.88
$ws.Range("Y5").Value = 11.17
$ws.Range("Z5").Value = 9.880000000000001
$ws.Range("AA5").Value = 19.76
$ws.Range("AB5").Value = 52277.81
$ws.Range("AC5").Value = 3898
$ws.Range("AD5").Value = 48.36
$ws.Range("AE5").Value = 50935
$ws.Range("AF5").Value = 3.7
$ws.Range("AG5").Value = 360
$ws.Range("AH5").Value = 0.19
$ws.Range("AI5").Value = 9.880000000000001
$ws.Range("AJ5").Value = 85026385

# Row 6
$ws.Range("D6").Value = 20213
$ws.Range("E6").Value = 2417
$ws.Range("F6").Value = 2417
$ws.Range("G6").Value = 2904
$ws.Range("H6").Value = 2149
$ws.Range("I6").Value = 1896
$ws.Range("K6").Value = 54330
$ws.Range("L6").Value = 8815
$ws.Range("M6").Value = 45515
$ws.Range("N6").Value = 43744
$ws.Range("P6").Value = 85
$ws.Range("Q6").Value = 2353
$ws.Range("R6").Value = -5053
$ws.Range("S6").Value = -992
$ws.Range("T6").Value = 230
$ws.Range("U6").Value = 2123
$ws.Range("V6").Value = 936
$ws.Range("W6").Value = 11.96
$ws.Range("X6").Value = 10.63
$ws.Range("Y6").Value = 4.36
$ws.Range("Z6").Value = 3.99
$ws.Range("AA6").Value = 19.37
$ws.Range("AB6").Value = 54071.11
$ws.Range("AC6").Value = 2226
$ws.Range("AD6").Value = 50.08
$ws.Range("AE6").Value = 52199
$ws.Range("AF6").Value = 2.14
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 85265865

# Row 7
$ws.Range("D7").Value = 22338
$ws.Range("E7").Value = 2232
$ws.Range("G7").Value = 2866
$ws.Range("H7").Value = 2214
$ws.Range("I7").Value = 2079
$ws.Range("K7").Value = 58947
$ws.Range("L7").Value = 10386
$ws.Range("M7").Value = 48562
$ws.Range("N7").Value = 46899
$ws.Range("P7").Value = 87
$ws.Range("Q7").Value = 2887
$ws.Range("R7").Value = -430
$ws.Range("S7").Value = 235
$ws.Range("T7").Value = 260
$ws.Range("U7").Value = 2546
$ws.Range("W7").Value = 9.99
$ws.Range("X7").Value = 9.91
$ws.Range("Y7").Value = 4.59
$ws.Range("Z7").Value = 3.91
$ws.Range("AA7").Value = 21.39
$ws.Range("AC7").Value = 2430
$ws.Range("AD7").Value = 36.87
$ws.Range("AE7").Value = 57386
$ws.Range("AF7").Value = 1.56
$ws.Range("AG7").Value = 50
$ws.Range("AH7").Value = 0.06
$ws.Range("AI7").Value = 2.06

# Row 8
$ws.Range("D8").Value = 26087
$ws.Range("E8").Value = 3532
$ws.Range("G8").Value = 4217
$ws.Range("H8").Value = 3275
$ws.Range("I8").Value = 3023
$ws.Range("K8").Value = 63096
$ws.Range("L8").Value = 11145
$ws.Range("M8").Value = 51951
$ws.Range("N8").Value = 50182
$ws.Range("P8").Value = 87
$ws.Range("Q8").Value = 3948
$ws.Range("R8").Value = -1724
$ws.Range("S8").Value = -100
$ws.Range("T8").Value = 237
$ws.Range("U8").Value = 3788
$ws.Range("W8").Value = 13.54
$ws.Range("X8").Value = 12.55
$ws.Range("Y8").Value = 6.23
$ws.Range("Z8").Value = 5.37
$ws.Range("AA8").Value = 21.45
$ws.Range("AC8").Value = 3525
$ws.Range("AD8").Value = 25.42
$ws.Range("AE8").Value = 61401
$ws.Range("AF8").Value = 1.46
$ws.Range("AG8").Value = 62
$ws.Range("AH8").Value = 0.07000000000000001
$ws.Range("AI8").Value = 1.77

# Row 9
$ws.Range("D9").Value = 27763
$ws.Range("E9").Value = 3901
$ws.Range("G9").Value = 4701
$ws.Range("H9").Value = 3633
$ws.Range("I9").Value = 3390
$ws.Range("K9").Value = 67744
$ws.Range("L9").Value = 11674
$ws.Range("M9").Value = 56070
$ws.Range("N9").Value = 53767
$ws.Range("P9").Value = 87
$ws.Range("Q9").Value = 4502
$ws.Range("R9").Value = -1457
$ws.Range("S9").Value = -243
$ws.Range("T9").Value = 262
$ws.Range("U9").Value = 4242
$ws.Range("W9").Value = 14.05
$ws.Range("X9").Value = 13.09
$ws.Range("Y9").Value = 6.52
$ws.Range("Z9").Value = 5.55
$ws.Range("AA9").Value = 20.82
$ws.Range("AC9").Value = 3954
$ws.Range("AD9").Value = 22.66
$ws.Range("AE9").Value = 65788
$ws.Range("AF9").Value = 1.36
$ws.Range("AG9").Value = 73
$ws.Range("AH9").Value = 0.08
$ws.Range("AI9").Value = 1.84
